$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# shift rows 93..113 down to 94..114
for ($r = 113; $r -ge 93; $r--) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r + 1, 1).Value = $a
    $ws.Cells.Item($r + 1, 2).Value = $b
    $ws.Cells.Item($r + 1, 3).Value = $c
    $ws.Cells.Item($r + 1, 4).Value = $d
}

# Now extend with new row 115 by copying format of row 114's A cell
$a114 = $ws.Cells.Item(114, 1).Value2
$ws.Cells.Item(115, 1).Value = $a114
$ws.Cells.Item(115, 1).Value = 44257
